$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maquinas")

# Insert a new row for the "Iberica" die-cutting machine above the existing
# "Descartonado" rows (row 15), shifting everything below it down by one.
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = "Troquelado"
$ws.Cells.Item(15, 2).Value = "Iberica"
$ws.Cells.Item(15, 3).Value = 2600
$ws.Cells.Item(15, 4).Value = 45
$ws.Cells.Item(15, 5).Value = 15

# Keep the "Observaciones" column readable after the insert.
$ws.Columns.Item(6).AutoFit()

# Reflect where the author ended up: looking at the newly added row on the
# "Maquinas" sheet.
$ws.Activate()
$ws.Range("E15").Select()
